{"js": "const replacements = [\n  [\"2024-07-06 Saturday\", \"2024-07-07 Sunday\"],\n  [\"835\u00d79=\", \"187\u00d77=\"],\n  [\"822\u00d75=\", \"465\u00d73=\"],\n  [\"340\u00d76=\", \"845\u00d74=\"],\n  [\"216\u00d76=\", \"711\u00d76=\"],\n  [\"369\u00d76=\", \"256\u00d74=\"],\n  [\"146\u00d78=\", \"251\u00d78=\"],\n  [\"749\u00d79=\", \"792\u00d74=\"],\n  [\"971\u00d72=\", \"982\u00d72=\"],\n  [\"556\u00d76=\", \"703\u00d77=\"],\n  [\"474\u00d76=\", \"577\u00d78=\"],\n  [\"751\u00d73=\", \"790\u00d75=\"],\n  [\"411\u00d73=\", \"277\u00d73=\"],\n  [\"386\u00d73=\", \"297\u00d79=\"],\n  [\"248\u00d79=\", \"922\u00d79=\"],\n  [\"357\u00d76=\", \"185\u00d74=\"],\n  [\"601\u00d76=\", \"345\u00d79=\"],\n  [\"823\u00d73=\", \"865\u00d75=\"],\n  [\"975\u00d79=\", \"689\u00d72=\"],\n  [\"764\u00d73=\", \"412\u00d77=\"],\n  [\"781\u00d73=\", \"862\u00d77=\"],\n  [\"439\u00d72=\", \"488\u00d75=\"],\n  [\"678\u00d74=\", \"586\u00d72=\"],\n  [\"120\u00d78=\", \"744\u00d76=\"],\n  [\"602\u00d73=\", \"662\u00d76=\"],\n  [\"713\u00d76=\", \"280\u00d77=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-06 Saturday\", \"2024-07-07 Sunday\"),\n    @(\"835\u00d79=\", \"187\u00d77=\"),\n    @(\"822\u00d75=\", \"465\u00d73=\"),\n    @(\"340\u00d76=\", \"845\u00d74=\"),\n    @(\"216\u00d76=\", \"711\u00d76=\"),\n    @(\"369\u00d76=\", \"256\u00d74=\"),\n    @(\"146\u00d78=\", \"251\u00d78=\"),\n    @(\"749\u00d79=\", \"792\u00d74=\"),\n    @(\"971\u00d72=\", \"982\u00d72=\"),\n    @(\"556\u00d76=\", \"703\u00d77=\"),\n    @(\"474\u00d76=\", \"577\u00d78=\"),\n    @(\"751\u00d73=\", \"790\u00d75=\"),\n    @(\"411\u00d73=\", \"277\u00d73=\"),\n    @(\"386\u00d73=\", \"297\u00d79=\"),\n    @(\"248\u00d79=\", \"922\u00d79=\"),\n    @(\"357\u00d76=\", \"185\u00d74=\"),\n    @(\"601\u00d76=\", \"345\u00d79=\"),\n    @(\"823\u00d73=\", \"865\u00d75=\"),\n    @(\"975\u00d79=\", \"689\u00d72=\"),\n    @(\"764\u00d73=\", \"412\u00d77=\"),\n    @(\"781\u00d73=\", \"862\u00d77=\"),\n    @(\"439\u00d72=\", \"488\u00d75=\"),\n    @(\"678\u00d74=\", \"586\u00d72=\"),\n    @(\"120\u00d78=\", \"744\u00d76=\"),\n    @(\"602\u00d73=\", \"662\u00d76=\"),\n    @(\"713\u00d76=\", \"280\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}\n"}
